$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Step 1: drop the rows that represent the "reference" category of each
# binary variable (the row whose label cell is blank, i.e. the 2nd row of a
# 2-row block). Delete from the highest index down so indices of the
# not-yet-processed rows stay stable.
$t.Rows(16).Delete()   # TVP.PREVIA = SIM  (old row 16, "SIM" / "43 (11.8)")
$t.Rows(14).Delete()   # SPT = SIM         (old row 14, "SIM" / "40 (11.0)")
$t.Rows(12).Delete()   # TEP = SIM         (old row 12, "SIM" / "7 ( 1.9)")
$t.Rows(10).Delete()   # Tipo.Atendimento = INTERNACAO (old row 10)
$t.Rows(4).Delete()    # Sexo = M          (old row 4, "M" / "141 (38.6)")

# --- Step 2: the middle "Categoria" column is no longer needed now that
# each remaining row only represents a single category -- remove it so the
# table goes from 3 columns back down to 2.
$t.Columns(2).Delete()

# --- Step 3: refresh the labels (fold the kept category into the bold
# label text) and the updated counts/statistics.
$d.Content.Find.Execute("Sexo (%)", $false, $false, $false, $false, $false, $true, 1, $false, "Sexo = M (%)", 2)
$d.Content.Find.Execute("Tipo.Atendimento (%)", $false, $false, $false, $false, $false, $true, 1, $false, "Tipo.Atendimento = INTERNAÇÃO (%)", 2)
$d.Content.Find.Execute("TEP (%)", $false, $false, $false, $false, $false, $true, 1, $false, "TEP = SIM (%)", 2)
$d.Content.Find.Execute("SPT (%)", $false, $false, $false, $false, $false, $true, 1, $false, "SPT = SIM (%)", 2)
$d.Content.Find.Execute("TVP.PREVIA (%)", $false, $false, $false, $false, $false, $true, 1, $false, "TVP.PREVIA = SIM (%)", 2)

$d.Content.Find.Execute("365", $false, $false, $false, $false, $false, $true, 1, $false, "426", 2)
$d.Content.Find.Execute("224 (61.4)", $false, $false, $false, $false, $false, $true, 1, $false, "159 (37.3)", 2)
$d.Content.Find.Execute("63.56 (15.97)", $false, $false, $false, $false, $false, $true, 1, $false, "63.43 (15.77)", 2)
$d.Content.Find.Execute("75.51 (19.47)", $false, $false, $false, $false, $false, $true, 1, $false, "76.29 (19.26)", 2)
$d.Content.Find.Execute("28.40 (6.20)", $false, $false, $false, $false, $false, $true, 1, $false, "28.71 (6.41)", 2)
$d.Content.Find.Execute("140 (38.4)", $false, $false, $false, $false, $false, $true, 1, $false, "250 (58.7)", 2)
$d.Content.Find.Execute("357 (98.1)", $false, $false, $false, $false, $false, $true, 1, $false, "8 ( 1.9)", 2)
$d.Content.Find.Execute("325 (89.0)", $false, $false, $false, $false, $false, $true, 1, $false, "54 (12.7)", 2)
$d.Content.Find.Execute("322 (88.2)", $false, $false, $false, $false, $false, $true, 1, $false, "62 (14.7)", 2)
